# Update the cryptos list sheet with refreshed price / volume data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$value)
    # Force the cell to keep a literal text value even when the new string
    # happens to look like a plain number (e.g. "557.30"), then drop the
    # temporary text number-format so the cell's style stays untouched.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "62.539.47"
$ws.Range("E2").Value = "  +4.34%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "2.427.11"
$ws.Range("E3").Value = "  +5.20%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.01%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "557.30"
$ws.Range("E5").Value = "  +2.98%  "

# Row 6 - Solana
Set-TextValue $ws.Range("D6") "138.91"
$ws.Range("E6").Value = "  +7.29%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.03%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  +2.03%  "

# Row 9 - LidoStakedEther
Set-TextValue $ws.Range("D9") "2.425.37"
$ws.Range("E9").Value = "  +5.24%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +3.72%  "

# Row 11 - Toncoin
Set-TextValue $ws.Range("D11") "5.78"
$ws.Range("E11").Value = "  +4.25%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  +0.47%  "

# Row 13 - Cardano
Set-TextValue $ws.Range("D13") "0.348"
$ws.Range("E13").Value = "  +4.74%  "

# Row 14 - Avalanche
Set-TextValue $ws.Range("D14") "26.15"
$ws.Range("E14").Value = "  +12.12%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D15") "2.860.21"
$ws.Range("E15").Value = "  +5.19%  "

# Row 16 - WrappedBTC
Set-TextValue $ws.Range("D16") "62.387.29"
$ws.Range("E16").Value = "  +4.11%  "

# Row 17 - ShibaInu
$ws.Range("E17").Value = "  +7.19%  "

# Row 18 - WrappedEther
Set-TextValue $ws.Range("D18") "2.427.07"
$ws.Range("E18").Value = "  +4.58%  "

# Row 19 - Chainlink
Set-TextValue $ws.Range("D19") "11.22"
$ws.Range("E19").Value = "  +7.02%  "

# Row 20 - BitcoinCash
Set-TextValue $ws.Range("D20") "345.57"
$ws.Range("E20").Value = "  +10.76%  "

# Row 21 - Polkadot
$ws.Range("E21").Value = "  +3.17%  "

# Row 22 - Uniswap
Set-TextValue $ws.Range("D22") "6.85"
$ws.Range("E22").Value = "  +4.36%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  +0.03%  "

# Row 24 - LEO
Set-TextValue $ws.Range("D24") "5.54"
$ws.Range("E24").Value = "  -2.83%  "

# Row 25 - Litecoin
Set-TextValue $ws.Range("D25") "65.27"
$ws.Range("E25").Value = "  +2.53%  "

# Row 26 - Kaspa
$ws.Range("E26").Value = "  +1.56%  "

# Row 27 - was Binance-PegBSC-USD, now Fetch.AI
$ws.Range("B27").Value = "Fetch.AI"
$ws.Range("C27").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue $ws.Range("D27") "1.56"
$ws.Range("E27").Value = "  +15.89%  "

# Row 28 - was Fetch.AI, now Binance-PegBSC-USD
$ws.Range("B28").Value = "Binance-PegBSC-USD"
$ws.Range("C28").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue $ws.Range("D28") "1.00"
$ws.Range("E28").Value = "  +0.04%  "

# Row 29 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D29") "8.15"
$ws.Range("E29").Value = "  +5.50%  "

# Row 30 - SuiNetwork
$ws.Range("E30").Value = "  +15.66%  "

# Row 31 - PancakeSwap
Set-TextValue $ws.Range("D31") "1.82"
$ws.Range("E31").Value = "  +5.54%  "

# Row 32 - PEPE
Set-TextValue $ws.Range("D32") "0.0₃0785"
$ws.Range("E32").Value = "  +8.19%  "

# Row 33 - Aptos
Set-TextValue $ws.Range("D33") "6.46"
$ws.Range("E33").Value = "  +10.94%  "

# Row 34 - Monero
Set-TextValue $ws.Range("D34") "172.30"
$ws.Range("E34").Value = "  +0.48%  "

# Row 35 - ImmutableX
$ws.Range("E35").Value = "  +6.02%  "

# Row 36 - PolygonEcosystemToken
Set-TextValue $ws.Range("D36") "0.397"
$ws.Range("E36").Value = "  +4.86%  "

# Row 37 - Bittensor
Set-TextValue $ws.Range("D37") "378.59"
$ws.Range("E37").Value = "  +19.44%  "

# Row 38 - EthereumClassic
Set-TextValue $ws.Range("D38") "18.56"
$ws.Range("E38").Value = "  +4.99%  "

# Row 39 - NEARProtocol
Set-TextValue $ws.Range("D39") "4.45"
$ws.Range("E39").Value = "  +11.39%  "

# Row 40 - USDe
$ws.Range("E40").Value = "  -0.01%  "

# Row 41 - FirstDigitalUSD
$ws.Range("E41").Value = "  -0.03%  "

# Row 42 - Stacks
Set-TextValue $ws.Range("D42") "1.70"
$ws.Range("E42").Value = "  +12.64%  "

# Row 43 - OKB
$ws.Range("E43").Value = "  +3.30%  "

# Row 44 - Aave
Set-TextValue $ws.Range("D44") "144.88"
$ws.Range("E44").Value = "  +6.51%  "

# Row 45 - Filecoin
$ws.Range("E45").Value = "  +7.35%  "

# Row 46 - InjectiveProtocol
Set-TextValue $ws.Range("D46") "20.81"
$ws.Range("E46").Value = "  +10.65%  "

# Row 47 - Mantle
Set-TextValue $ws.Range("D47") "0.593"
$ws.Range("E47").Value = "  +4.19%  "

# Row 48 - Stellar
Set-TextValue $ws.Range("D48") "0.0952"
$ws.Range("E48").Value = "  +1.65%  "

# Row 49 - Hedera
$ws.Range("E49").Value = "  +6.06%  "

# Row 50 - VeChain
$ws.Range("E50").Value = "  +4.55%  "

# Row 51 - EnergySwap
Set-TextValue $ws.Range("D51") "17.85"
$ws.Range("E51").Value = "  +6.38%  "
